$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update existing rows 174-176 with revised figures
$ws.Range("B174").Value = -1033
$ws.Range("C174").Value = -1015
$ws.Range("D174").Value = 1212
$ws.Range("I174").Value = 3477
$ws.Range("J174").Value = 3439
$ws.Range("M174").Value = 1447
$ws.Range("N174").Value = 889
$ws.Range("O174").Value = 580
$ws.Range("P174").Value = 250

$ws.Range("B175").Value = -816
$ws.Range("C175").Value = -796
$ws.Range("D175").Value = 1287
$ws.Range("I175").Value = 3406
$ws.Range("J175").Value = 3370
$ws.Range("M175").Value = 1410
$ws.Range("N175").Value = 792
$ws.Range("O175").Value = 599
$ws.Range("P175").Value = 260

$ws.Range("B176").Value = -835
$ws.Range("C176").Value = -800
$ws.Range("D176").Value = 1322
$ws.Range("E176").Value = 2122
$ws.Range("I176").Value = 3496
$ws.Range("J176").Value = 3445
$ws.Range("L176").Value = 81
$ws.Range("M176").Value = 1503
$ws.Range("N176").Value = 814
$ws.Range("O176").Value = 581
$ws.Range("P176").Value = 269

# Add new row 177 with the August 2021 data
# Force the date-looking text to stay a text string (not get auto-converted
# to a date serial) by setting the format to Text before entry, then clear
# the formatting afterward so the cell stays unstyled like the rest of column A.
$ws.Range("A177").NumberFormat = "@"
$ws.Range("A177").Value = "01-08-2021"
$ws.Range("A177").ClearFormats()

$ws.Range("B177").Value = -1008
$ws.Range("C177").Value = -974
$ws.Range("D177").Value = 1255
$ws.Range("E177").Value = 2229
$ws.Range("F177").Value = -35
$ws.Range("G177").Value = 8
$ws.Range("H177").Value = 43
$ws.Range("I177").Value = 3535
$ws.Range("J177").Value = 3485
$ws.Range("K177").Value = 63
$ws.Range("L177").Value = 57
$ws.Range("M177").Value = 1526
$ws.Range("N177").Value = 903
$ws.Range("O177").Value = 536
$ws.Range("P177").Value = 263
$ws.Range("Q177").Value = 136
$ws.Range("R177").Value = 51
